$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.437724
$ws.Range("H2").Value = 1.313172
$ws.Range("I2").Value = 0.02046276855287852
$ws.Range("J2").Value = 0.02204588088728605
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.506715
$ws.Range("N2").Value = 1.520145
$ws.Range("O2").Value = 0.003122343715987576
$ws.Range("P2").Value = 0.003132472094339857
$ws.Range("Q2").Value = 0.22180131666
$ws.Range("R2").Value = 1.99621184994
$ws.Range("S2").Value = [double]"6.389179680278845E-05"
$ws.Range("T2").Value = [double]"6.905810667456395E-05"

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.437724
$ws.Range("H3").Value = 1.313172
$ws.Range("I3").Value = 0.02046276855287852
$ws.Range("J3").Value = 0.02204588088728605
$ws.Range("M3").Value = 88.13219433333332
$ws.Range("N3").Value = 264.396583
$ws.Range("O3").Value = 0.5430646480820168
$ws.Range("P3").Value = 0.5448262620252092
$ws.Range("Q3").Value = 38.57757663236399
$ws.Range("R3").Value = 347.198189691276
$ws.Range("S3").Value = 0.01111260620295274
$ws.Range("T3").Value = 0.01201117487687306

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.437724
$ws.Range("H4").Value = 1.313172
$ws.Range("I4").Value = 0.02046276855287852
$ws.Range("J4").Value = 0.02204588088728605
$ws.Range("M4").Value = 1.5741895
$ws.Range("N4").Value = 3.148379
$ws.Range("O4").Value = 0.009700049718478087
$ws.Range("P4").Value = 0.006487676741301404
$ws.Range("Q4").Value = 0.6890605246979999
$ws.Range("R4").Value = 4.134363148188
$ws.Range("S4").Value = 0.0001984898723406316
$ws.Range("T4").Value = 0.0001430265486739469

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.437724
$ws.Range("H5").Value = 1.313172
$ws.Range("I5").Value = 0.02046276855287852
$ws.Range("J5").Value = 0.02204588088728605
$ws.Range("M5").Value = 72.07364666666666
$ws.Range("N5").Value = 216.22094
$ws.Range("O5").Value = 0.4441129584835175
$ws.Range("P5").Value = 0.4455535891391496
$ws.Range("Q5").Value = 31.54836491352
$ws.Range("R5").Value = 283.93528422168
$ws.Range("S5").Value = 0.009087780680782366
$ws.Range("T5").Value = 0.009822621355064478

$ws.Range("G6").Value = 12.48419333333333
$ws.Range("I6").Value = 0.5836124104444559
$ws.Range("J6").Value = 0.6287638767819841
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.506715
$ws.Range("N6").Value = 1.520145
$ws.Range("O6").Value = 0.003122343715987576
$ws.Range("P6").Value = 0.003132472094339857
$ws.Range("Q6").Value = 6.3259280249
$ws.Range("R6").Value = 56.9333522241
$ws.Range("S6").Value = 0.001822238542323609
$ws.Range("T6").Value = 0.001969585297948509

$ws.Range("G7").Value = 12.48419333333333
$ws.Range("I7").Value = 0.5836124104444559
$ws.Range("J7").Value = 0.6287638767819841
$ws.Range("M7").Value = 88.13219433333332
$ws.Range("N7").Value = 264.396583
$ws.Range("O7").Value = 0.5430646480820168
$ws.Range("P7").Value = 0.5448262620252092
$ws.Range("Q7").Value = 1100.259352948237
$ws.Range("R7").Value = 9902.334176534137
$ws.Range("S7").Value = 0.316939268294316
$ws.Range("T7").Value = 0.3425670726836076

$ws.Range("G8").Value = 12.48419333333333
$ws.Range("I8").Value = 0.5836124104444559
$ws.Range("J8").Value = 0.6287638767819841
$ws.Range("M8").Value = 1.5741895
$ws.Range("N8").Value = 3.148379
$ws.Range("O8").Value = 0.009700049718478087
$ws.Range("P8").Value = 0.006487676741301404
$ws.Range("Q8").Value = 19.65248606130333
$ws.Range("R8").Value = 117.91491636782
$ws.Range("S8").Value = 0.005661069397632062
$ws.Range("T8").Value = 0.00407921677916898

$ws.Range("G9").Value = 12.48419333333333
$ws.Range("I9").Value = 0.5836124104444559
$ws.Range("J9").Value = 0.6287638767819841
$ws.Range("M9").Value = 72.07364666666666
$ws.Range("N9").Value = 216.22094
$ws.Range("O9").Value = 0.4441129584835175
$ws.Range("P9").Value = 0.4455535891391496
$ws.Range("Q9").Value = 899.7813392250221
$ws.Range("R9").Value = 8098.032053025199
$ws.Range("S9").Value = 0.2591898342101842
$ws.Range("T9").Value = 0.280148002021259

$ws.Range("G10").Value = 1.796802333333333
$ws.Range("H10").Value = 5.390407
$ws.Range("I10").Value = 0.08399710841140098
$ws.Range("J10").Value = 0.09049558675938332
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.506715
$ws.Range("N10").Value = 1.520145
$ws.Range("O10").Value = 0.003122343715987576
$ws.Range("P10").Value = 0.003132472094339857
$ws.Range("Q10").Value = 0.9104666943349999
$ws.Range("R10").Value = 8.194200249014999
$ws.Range("S10").Value = 0.0002622678436094651
$ws.Range("T10").Value = 0.0002834749001846797

$ws.Range("G11").Value = 1.796802333333333
$ws.Range("H11").Value = 5.390407
$ws.Range("I11").Value = 0.08399710841140098
$ws.Range("J11").Value = 0.09049558675938332
$ws.Range("M11").Value = 88.13219433333332
$ws.Range("N11").Value = 264.396583
$ws.Range("O11").Value = 0.5430646480820168
$ws.Range("P11").Value = 0.5448262620252092
$ws.Range("Q11").Value = 158.3561324199201
$ws.Range("R11").Value = 1425.205191779281
$ws.Range("S11").Value = 0.04561586011934449
$ws.Range("T11").Value = 0.04930437226389283

$ws.Range("G12").Value = 1.796802333333333
$ws.Range("H12").Value = 5.390407
$ws.Range("I12").Value = 0.08399710841140098
$ws.Range("J12").Value = 0.09049558675938332
$ws.Range("M12").Value = 1.5741895
$ws.Range("N12").Value = 3.148379
$ws.Range("O12").Value = 0.009700049718478087
$ws.Range("P12").Value = 0.006487676741301404
$ws.Range("Q12").Value = 2.828507366708833
$ws.Range("R12").Value = 16.971044200253
$ws.Range("S12").Value = 0.0008147761277989834
$ws.Range("T12").Value = 0.0005871061134092745

$ws.Range("G13").Value = 1.796802333333333
$ws.Range("H13").Value = 5.390407
$ws.Range("I13").Value = 0.08399710841140098
$ws.Range("J13").Value = 0.09049558675938332
$ws.Range("M13").Value = 72.07364666666666
$ws.Range("N13").Value = 216.22094
$ws.Range("O13").Value = 0.4441129584835175
$ws.Range("P13").Value = 0.4455535891391496
$ws.Range("Q13").Value = 129.5020965025089
$ws.Range("R13").Value = 1165.51886852258
$ws.Range("S13").Value = 0.03730420432064804
$ws.Range("T13").Value = 0.04032063348189654

$ws.Range("G14").Value = 4.608308
$ws.Range("H14").Value = 9.216616
$ws.Range("I14").Value = 0.2154296772038511
$ws.Range("J14").Value = 0.154731001361478
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.506715
$ws.Range("N14").Value = 1.520145
$ws.Range("O14").Value = 0.003122343715987576
$ws.Range("P14").Value = 0.003132472094339857
$ws.Range("Q14").Value = 2.33509878822
$ws.Range("R14").Value = 14.01059272932
$ws.Range("S14").Value = 0.0006726454988546766
$ws.Range("T14").Value = 0.0004846905438940922

$ws.Range("G15").Value = 4.608308
$ws.Range("H15").Value = 9.216616
$ws.Range("I15").Value = 0.2154296772038511
$ws.Range("J15").Value = 0.154731001361478
$ws.Range("M15").Value = 88.13219433333332
$ws.Range("N15").Value = 264.396583
$ws.Range("O15").Value = 0.5430646480820168
$ws.Range("P15").Value = 0.5448262620252092
$ws.Range("Q15").Value = 406.1402962038546
$ws.Range("R15").Value = 2436.841777223128
$ws.Range("S15").Value = 0.1169922418371319
$ws.Range("T15").Value = 0.08430151309119162

$ws.Range("G16").Value = 4.608308
$ws.Range("H16").Value = 9.216616
$ws.Range("I16").Value = 0.2154296772038511
$ws.Range("J16").Value = 0.154731001361478
$ws.Range("M16").Value = 1.5741895
$ws.Range("N16").Value = 3.148379
$ws.Range("O16").Value = 0.009700049718478087
$ws.Range("P16").Value = 0.006487676741301404
$ws.Range("Q16").Value = 7.254350066365999
$ws.Range("R16").Value = 29.017400265464
$ws.Range("S16").Value = 0.002089678579713041
$ws.Range("T16").Value = 0.001003844718691137

$ws.Range("G17").Value = 4.608308
$ws.Range("H17").Value = 9.216616
$ws.Range("I17").Value = 0.2154296772038511
$ws.Range("J17").Value = 0.154731001361478
$ws.Range("M17").Value = 72.07364666666666
$ws.Range("N17").Value = 216.22094
$ws.Range("O17").Value = 0.4441129584835175
$ws.Range("P17").Value = 0.4455535891391496
$ws.Range("Q17").Value = 332.1375625231733
$ws.Range("R17").Value = 1992.82537513904
$ws.Range("S17").Value = 0.0956751112881515
$ws.Range("T17").Value = 0.06894095300770116

$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 2.064212666666667
$ws.Range("H18").Value = 6.192638000000001
$ws.Range("I18").Value = 0.09649803538741349
$ws.Range("J18").Value = 0.1039636542098684
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.506715
$ws.Range("N18").Value = 1.520145
$ws.Range("O18").Value = 0.003122343715987576
$ws.Range("P18").Value = 0.003132472094339857
$ws.Range("Q18").Value = 1.04596752139
$ws.Range("R18").Value = 9.413707692510002
$ws.Range("S18").Value = 0.0003013000343970373
$ws.Range("T18").Value = 0.0003256632456380111

$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 2.064212666666667
$ws.Range("H19").Value = 6.192638000000001
$ws.Range("I19").Value = 0.09649803538741349
$ws.Range("J19").Value = 0.1039636542098684
$ws.Range("M19").Value = 88.13219433333332
$ws.Range("N19").Value = 264.396583
$ws.Range("O19").Value = 0.5430646480820168
$ws.Range("P19").Value = 0.5448262620252092
$ws.Range("Q19").Value = 181.9235918839948
$ws.Range("R19").Value = 1637.312326955954
$ws.Range("S19").Value = 0.05240467162827171
$ws.Range("T19").Value = 0.05664212910964401

$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 0.6666666666666666
$ws.Range("G20").Value = 2.064212666666667
$ws.Range("H20").Value = 6.192638000000001
$ws.Range("I20").Value = 0.09649803538741349
$ws.Range("J20").Value = 0.1039636542098684
$ws.Range("M20").Value = 1.5741895
$ws.Range("N20").Value = 3.148379
$ws.Range("O20").Value = 0.009700049718478087
$ws.Range("P20").Value = 0.006487676741301404
$ws.Range("Q20").Value = 3.249461905633666
$ws.Range("R20").Value = 19.496771433802
$ws.Range("S20").Value = 0.0009360357409933687
$ws.Range("T20").Value = 0.0006744825813580651

$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 0.6666666666666666
$ws.Range("G21").Value = 2.064212666666667
$ws.Range("H21").Value = 6.192638000000001
$ws.Range("I21").Value = 0.09649803538741349
$ws.Range("J21").Value = 0.1039636542098684
$ws.Range("M21").Value = 72.07364666666666
$ws.Range("N21").Value = 216.22094
$ws.Range("O21").Value = 0.4441129584835175
$ws.Range("P21").Value = 0.4455535891391496
$ws.Range("Q21").Value = 148.7753343821911
$ws.Range("R21").Value = 1338.97800943972
$ws.Range("S21").Value = 0.04285602798375136
$ws.Range("T21").Value = 0.04632137927322832
